$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 7.186737775649103

# Row 3
$ws.Range("D3").Value = 19
$ws.Range("E3").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("F3").Value = "Jewel"
$ws.Range("H3").Value = 6.491841461046875
$ws.Range("I3").Value = "Black or African American"

# Row 4
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("F4").Value = "Colleen"
$ws.Range("H4").Value = 6.098775820819761
$ws.Range("I4").Value = "White"

# Row 5
$ws.Range("H5").Value = 5.03843188797754

# Row 6
$ws.Range("H6").Value = 5.020761003118488

# Row 7
$ws.Range("H7").Value = 4.239120605821088

# Row 8
$ws.Range("H8").Value = 1.248319253184411

# Row 9
$ws.Range("H9").Value = 1.073514296423548

# Row 10
$ws.Range("D10").Value = 32
$ws.Range("E10").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("F10").Value = "Kellie"
$ws.Range("H10").Value = 0.4217658868825644
$ws.Range("I10").Value = "White"

# Row 11
$ws.Range("D11").Value = 21
$ws.Range("E11").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("F11").Value = "Bri"
$ws.Range("H11").Value = 0.4004490700212808

# Row 12
$ws.Range("D12").Value = 30
$ws.Range("E12").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("F12").Value = "Shadaisia"
$ws.Range("H12").Value = 0.3953773248513697
$ws.Range("I12").Value = "Black or African American"

# Row 13
$ws.Range("D13").Value = 33
$ws.Range("E13").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("F13").Value = "Shaniek"
$ws.Range("H13").Value = 0.1805031549777598

# Row 14
$ws.Range("H14").Value = 13.04325402792447

# Row 15
$ws.Range("H15").Value = 8.342377812971202

# Row 16
$ws.Range("H16").Value = 7.489472321657063

# Row 17
$ws.Range("H17").Value = 7.219432926815826

# Row 18
$ws.Range("D18").Value = 22
$ws.Range("E18").Value = "60db4fde6193c50664c9c478"
$ws.Range("F18").Value = "Edosagbe"
$ws.Range("H18").Value = 5.404956080902719
$ws.Range("I18").Value = "Black or African American"

# Row 19
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("F19").Value = "Jamarii"
$ws.Range("H19").Value = 5.203546488046102

# Row 20
$ws.Range("D20").Value = 26
$ws.Range("E20").Value = "5dd671942b033b5ec8bc97b4"
$ws.Range("F20").Value = "Juan"
$ws.Range("H20").Value = 5.194694186643499
$ws.Range("I20").Value = "Hispanic"

# Row 21
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = "60b322994d0b901954690036"
$ws.Range("F21").Value = "Brennan"
$ws.Range("H21").Value = 4.334666484926464

# Row 22
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = "5e2522d6b734b47915f88275"
$ws.Range("F22").Value = "Corey"
$ws.Range("H22").Value = 4.178693876440433

# Row 23
$ws.Range("H23").Value = 3.419194189605884

# Row 24
$ws.Range("H24").Value = 2.385885516067507

# Row 25
$ws.Range("H25").Value = 2.223286854337817
